$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Replace row 8 data: "Caldas" (mislabeled under Antioquia) -> "Manizales" (Caldas dept)
$ws.Range("A8").Value = 17001
$ws.Range("B8").Value = "Manizales"
$ws.Range("C8").Value = 17
$ws.Range("D8").Value = "Caldas"
$ws.Range("E8").Value = 5.07
$ws.Range("F8").Value = -75.5205556

# Update the active selection to A9, matching the saved workbook view state
$ws.Range("A9").Select()
